# escena.xlsx - "Creande escena y parametros mejorados"
# Updates the mean (B column) values of the existing Gaussian-scene rows
# and appends new rows (23-40) extending the parameter grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Update the "Medias" (B column) of the existing rows ---
$ws.Range("B3").Value = 2.5
$ws.Range("B4").Value = 2.5
$ws.Range("B9").Value = 2.5
$ws.Range("B10").Value = 2.5
$ws.Range("B11").Value = -2.5
$ws.Range("B12").Value = -2.5
$ws.Range("B13").Value = -2.5
$ws.Range("B14").Value = -2.5
$ws.Range("B15").Value = -7.5
$ws.Range("B16").Value = -7.5
$ws.Range("B17").Value = -7.5
$ws.Range("B18").Value = -7.5
$ws.Range("B19").Value = -12.5
$ws.Range("B20").Value = -12.5
$ws.Range("B21").Value = -12.5
$ws.Range("B22").Value = -12.5

# --- 2. Append new parameter rows 23-40 ---
# Copy the formatting (borders / centered alignment / fill) from an existing
# data row (row 9) onto each new row before writing its values.

$ws.Range("B9:F9").Copy()
$ws.Range("B23:F23").PasteSpecial(-4122)
$ws.Range("B23").Value = -2.5
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0.9
$ws.Range("E23").Value = 0.9
$ws.Range("F23").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B24:F24").PasteSpecial(-4122)
$ws.Range("B24").Value = 2.5
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.9
$ws.Range("E24").Value = 0.9
$ws.Range("F24").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B25:F25").PasteSpecial(-4122)
$ws.Range("B25").Value = -7.5
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0.9
$ws.Range("E25").Value = 0.9
$ws.Range("F25").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B26:F26").PasteSpecial(-4122)
$ws.Range("B26").Value = -12.5
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.9
$ws.Range("E26").Value = 0.9
$ws.Range("F26").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B27:F27").PasteSpecial(-4122)
$ws.Range("B27").Value = 7.5
$ws.Range("C27").Value = -10
$ws.Range("D27").Value = 0.9
$ws.Range("E27").Value = 0.9
$ws.Range("F27").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B28:F28").PasteSpecial(-4122)
$ws.Range("B28").Value = 7.5
$ws.Range("C28").Value = -5
$ws.Range("D28").Value = 0.9
$ws.Range("E28").Value = 0.9
$ws.Range("F28").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B29:F29").PasteSpecial(-4122)
$ws.Range("B29").Value = 7.5
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0.9
$ws.Range("E29").Value = 0.9
$ws.Range("F29").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B30:F30").PasteSpecial(-4122)
$ws.Range("B30").Value = 7.5
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 0.9
$ws.Range("E30").Value = 0.9
$ws.Range("F30").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B31:F31").PasteSpecial(-4122)
$ws.Range("B31").Value = 7.5
$ws.Range("C31").Value = 10
$ws.Range("D31").Value = 0.9
$ws.Range("E31").Value = 0.9
$ws.Range("F31").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B32:F32").PasteSpecial(-4122)
$ws.Range("B32").Value = 12.5
$ws.Range("C32").Value = -10
$ws.Range("D32").Value = 0.9
$ws.Range("E32").Value = 0.9
$ws.Range("F32").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B33:F33").PasteSpecial(-4122)
$ws.Range("B33").Value = 12.5
$ws.Range("C33").Value = -5
$ws.Range("D33").Value = 0.9
$ws.Range("E33").Value = 0.9
$ws.Range("F33").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B34:F34").PasteSpecial(-4122)
$ws.Range("B34").Value = 12.5
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0.9
$ws.Range("E34").Value = 0.9
$ws.Range("F34").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B35:F35").PasteSpecial(-4122)
$ws.Range("B35").Value = 12.5
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = 0.9
$ws.Range("E35").Value = 0.9
$ws.Range("F35").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B36:F36").PasteSpecial(-4122)
$ws.Range("B36").Value = 12.5
$ws.Range("C36").Value = 10
$ws.Range("D36").Value = 0.9
$ws.Range("E36").Value = 0.9
$ws.Range("F36").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B37:F37").PasteSpecial(-4122)
$ws.Range("B37").Value = 7.5
$ws.Range("C37").Value = -15
$ws.Range("D37").Value = 0.9
$ws.Range("E37").Value = 0.9
$ws.Range("F37").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B38:F38").PasteSpecial(-4122)
$ws.Range("B38").Value = 2.5
$ws.Range("C38").Value = -15
$ws.Range("D38").Value = 0.9
$ws.Range("E38").Value = 0.9
$ws.Range("F38").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B39:F39").PasteSpecial(-4122)
$ws.Range("B39").Value = -2.5
$ws.Range("C39").Value = -15
$ws.Range("D39").Value = 0.9
$ws.Range("E39").Value = 0.9
$ws.Range("F39").Value = 40000

$ws.Range("B9:F9").Copy()
$ws.Range("B40:F40").PasteSpecial(-4122)
$ws.Range("B40").Value = -7.5
$ws.Range("C40").Value = -15
$ws.Range("D40").Value = 0.9
$ws.Range("E40").Value = 0.9
$ws.Range("F40").Value = 40000

$excel.CutCopyMode = $false

# --- 3. Restore the scrolled view / active selection (I31) ---
[void]$ws.Range("A26").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I31").Select()
